$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$src = $ws.Cells.Item(118, 1)
$st = $src.Style
Write-Host "style name:" $st.Name
